$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Mint the "Heading 2" paragraph style (and its linked "Heading 2
#    Char" character style) so the new "Simulator Test:" heading can
#    use it. Word normally materializes the built-in Heading 2
#    definition into styles.xml the first time the style is actually
#    applied to a paragraph, so we do that against a throwaway
#    scratch paragraph, refine the style's properties to match the
#    canonical built-in template, then remove the scratch paragraph
#    again (the style definition itself stays behind in styles.xml).
# ------------------------------------------------------------------
$scratch = $d.Paragraphs.Last.Range.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "scratch"
$d.Paragraphs.Last.Style = "Heading 2"

$heading2 = $d.Styles("Heading 2")
$heading2.NameLocal = "heading 2"
$heading2.NextParagraphStyle = "Normal"
$heading2.LinkStyle = "Heading2Char"
$heading2.UnhideWhenUsed = $true
$heading2.QuickStyle = $true
$heading2.Priority = 9
$heading2.Font.Bold = $true
$heading2.Font.BoldBi = $true
$heading2.Font.Size = 13
$heading2.Font.SizeBi = 13
$heading2.Font.TextColor.ObjectThemeColor = 4
$heading2.ParagraphFormat.SpaceBefore = 10
$heading2.ParagraphFormat.SpaceAfter = 0

$heading2Char = $d.Styles.Add("Heading 2 Char", 2)
$heading2Char.BaseStyle = "DefaultParagraphFont"
$heading2Char.LinkStyle = "Heading2"
$heading2Char.Priority = 9
$heading2Char.Font.Bold = $true
$heading2Char.Font.BoldBi = $true
$heading2Char.Font.Size = 13
$heading2Char.Font.SizeBi = 13
$heading2Char.Font.TextColor.ObjectThemeColor = 4

# Remove the scratch paragraph (and its paragraph mark) again -- only
# its styling side-effect (the minted styles above) should remain.
$d.Paragraphs.Last.Range.Delete()

# ------------------------------------------------------------------
# 2. Drop the stray "_GoBack" bookmark that used to sit right after
#    "...ignore the heeling in this testing." -- it is being relocated
#    to the very end of the document (see step 3).
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 3. Append the new "Simulator Test" section at the end of the body:
#      - a blank spacer paragraph
#      - the "Simulator Test:" Heading 2 paragraph
#      - the paragraph describing the simulator test
#      - a trailing empty paragraph holding the (relocated) _GoBack
#        bookmark
#    Building this as one Open-XML fragment keeps the paragraph marks
#    clean (no stray empty runs) and lets us place the
#    lastRenderedPageBreak marker exactly where Word would.
# ------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)

$openXmlPackage = @"
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p/>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading2"/>
            </w:pPr>
            <w:r>
              <w:lastRenderedPageBreak/>
              <w:t>Simulator Test:</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:r>
              <w:t>In order to test tacking algorithm, options within simulation: n-parameter sweep, n-parameter perturbation+hill climbing, or genetic algorithm.  Parameters to refine: sail set values, which angles define certain behaviors.  Want to measure exit velocity, time spent in tack normalized by entrance velocity</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$endRange.InsertXML($openXmlPackage)
